$wb = $excel.ActiveWorkbook

# --- "survey" sheet: remove two rows that are no longer part of the in-person visit ---
# Row 26 (note / display_dates / "Display dates ...") is removed first so the
# row-24 deletion below doesn't shift it before we get to it.
$survey = $wb.Worksheets.Item("survey")
$survey.Rows.Item(26).Delete()
$survey.Rows.Item(24).Delete()

# --- "settings" sheet: rename the form's display title ---
$settings = $wb.Worksheets.Item("settings")
$settings.Range("A2").Value = "Form - In-person Visit"
